$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 71429864
$ws.Range("I8").Value = 125000020
$ws.Range("J8").Value = 3000
$ws.Range("K8").Value = 375000060
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = -374999921
$ws.Range("N8").Value = -9278

$ws.Range("H82").Value = 1268.4
$ws.Range("I82").Value = 585.5
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 1756.5
$ws.Range("L82").Value = 12000
$ws.Range("M82").Value = -1350.5
$ws.Range("N82").Value = -12812

$ws.Range("H85").Value = 1268.4
$ws.Range("I85").Value = 585.5
$ws.Range("J85").Value = 4000
$ws.Range("K85").Value = 1756.5
$ws.Range("L85").Value = 12000
$ws.Range("M85").Value = -352.5
$ws.Range("N85").Value = -14808

$ws.Range("H111").Value = 590.25
$ws.Range("I111").Value = 643
$ws.Range("J111").Value = 432
$ws.Range("K111").Value = 1929
$ws.Range("L111").Value = 1296
$ws.Range("M111").Value = 1138
$ws.Range("N111").Value = -7430

$ws.Range("I112").Value = 910
$ws.Range("J112").Value = 1994
$ws.Range("K112").Value = 2730
$ws.Range("L112").Value = 5982
$ws.Range("M112").Value = -1622
$ws.Range("N112").Value = -8198

$ws.Range("H124").Value = 39795
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 39795
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 39795
$ws.Range("N124").Value = -49615

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 295.36365
$ws.Range("I4").Value = 274.9
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 274.9
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -158.9
$ws.Range("N4").Value = -732

$ws.Range("H10").Value = 200
$ws.Range("I10").Value = 200
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 200
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -30

$ws.Range("H97").Value = 1109.2222
$ws.Range("I97").Value = 1117
$ws.Range("J97").Value = 1070.3334
$ws.Range("K97").Value = 1117
$ws.Range("L97").Value = 1070.3334
$ws.Range("M97").Value = -621
$ws.Range("N97").Value = -2062.3334

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 76.666664
$ws.Range("I13").Value = 76.666664
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 229.999992
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -61.99999199999999

$ws.Range("H58").Value = 2875.5557
$ws.Range("I58").Value = 1500
$ws.Range("J58").Value = 3047.5
$ws.Range("K58").Value = 4500
$ws.Range("L58").Value = 9142.5
$ws.Range("M58").Value = -4372
$ws.Range("N58").Value = -9398.5

$ws.Range("H64").Value = 3061
$ws.Range("I64").Value = 2432
$ws.Range("J64").Value = 3118.182
$ws.Range("K64").Value = 7296
$ws.Range("L64").Value = 9354.545999999998
$ws.Range("M64").Value = -7026
$ws.Range("N64").Value = -9894.545999999998

$ws.Range("H67").Value = 3061
$ws.Range("I67").Value = 2432
$ws.Range("J67").Value = 3118.182
$ws.Range("K67").Value = 7296
$ws.Range("L67").Value = 9354.545999999998
$ws.Range("M67").Value = -6360
$ws.Range("N67").Value = -11226.546

$ws.Range("H70").Value = 3844.6
$ws.Range("I70").Value = 911.5
$ws.Range("J70").Value = 5800
$ws.Range("K70").Value = 2734.5
$ws.Range("L70").Value = 17400
$ws.Range("M70").Value = -2419.5
$ws.Range("N70").Value = -18030

$ws.Range("H73").Value = 3844.6
$ws.Range("I73").Value = 911.5
$ws.Range("J73").Value = 5800
$ws.Range("K73").Value = 2734.5
$ws.Range("L73").Value = 17400
$ws.Range("M73").Value = -1642.5
$ws.Range("N73").Value = -19584

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()

$ws.Range("H103").Value = 756.25
$ws.Range("I103").Value = 341.66666
$ws.Range("J103").Value = 2000
$ws.Range("K103").Value = 1024.99998
$ws.Range("L103").Value = 6000
$ws.Range("M103").Value = -145.9999800000001
$ws.Range("N103").Value = -7758

$ws.Range("H108").Value = 2340.75
$ws.Range("I108").Value = 475.33334
$ws.Range("J108").Value = 3460
$ws.Range("K108").Value = 1426.00002
$ws.Range("L108").Value = 10380
$ws.Range("M108").Value = 1453.99998
$ws.Range("N108").Value = -16140

$ws.Range("H109").Value = 2191.375
$ws.Range("I109").Value = 1020.25
$ws.Range("J109").Value = 3362.5
$ws.Range("K109").Value = 3060.75
$ws.Range("L109").Value = 10087.5
$ws.Range("M109").Value = -2020.75
$ws.Range("N109").Value = -12167.5

$ws.Range("H113").Value = 552.3
$ws.Range("I113").Value = 540.2029
$ws.Range("J113").Value = 579.2258
$ws.Range("K113").Value = 1620.6087
$ws.Range("L113").Value = 1737.6774
$ws.Range("M113").Value = 549.3913
$ws.Range("N113").Value = -6077.6774

$ws.Range("H124").Value = 1047.125
$ws.Range("I124").Value = 529.5
$ws.Range("J124").Value = 2600
$ws.Range("K124").Value = 1588.5
$ws.Range("L124").Value = 7800
$ws.Range("M124").Value = 3321.5
$ws.Range("N124").Value = -17620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4121.25
$ws.Range("I80").Value = 4210
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 4210
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -3212
$ws.Range("N80").Value = -5496

$ws.Range("H83").Value = 4121.25
$ws.Range("I83").Value = 4210
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 21050
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -16058
$ws.Range("N83").Value = -27484

$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 675
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 700
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 700
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -1076

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 5000750
$ws.Range("I11").Value = 10000500
$ws.Range("J11").Value = 999.5
$ws.Range("K11").Value = 10000500
$ws.Range("L11").Value = 999.5
$ws.Range("M11").Value = -10000358
$ws.Range("N11").Value = -1283.5

$ws.Range("H19").Value = 7968.3335
$ws.Range("I19").Value = 8905
$ws.Range("J19").Value = 7500
$ws.Range("K19").Value = 8905
$ws.Range("L19").Value = 7500
$ws.Range("M19").Value = -8731
$ws.Range("N19").Value = -7848
